$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet gained a new (blank) column inserted right
# before the existing "Late" column (old column N), pushing the old
# N/O/P columns ("Late", blank, "Outstanding") one column to the right.
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Make this sheet the active one (matches workbook activeTab changing
# from the "Transactions" sheet to the "Repayment Schedule" sheet, and
# the "Repayment Schedule" sheetView gaining tabSelected="1" while
# "Transactions" loses it).
$ws.Activate()

# Insert a new blank column at N, shifting old N/O/P to O/P/Q.
$ws.Columns("N").Insert()

# Update the remembered selection on the sheet to match the new layout.
$ws.Range("S5").Select()
